# Additions to Forgot Password
# Adds a new "ForgotPasswordEmail" worksheet after "ForgotPasswordHeaderFooter"
# containing the default/validation text used by the forgot-password e-mail flow.

$wb = $excel.ActiveWorkbook

# Locate the sheet after which the new sheet must be inserted.
$afterSheet = $wb.Worksheets.Item("ForgotPasswordHeaderFooter")

# Add the new worksheet right after it.
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "ForgotPasswordEmail"

# Enter the header/value pairs column by column (matches the order the
# strings were originally authored in, so the shared-string table comes
# out in the same sequence).
$newSheet.Range("B1").Value = "E-Mail Address"
$newSheet.Range("B2").Value = "me.automaton@gmail.com"

$newSheet.Range("A1").Value = "Default E-Mail Address Text"
$newSheet.Range("A2").Value = "Enter e-mail address"

$newSheet.Range("C1").Value = "Missing E-Mail Address Text"
$newSheet.Range("C2").Value = "E-mail Address is Required."

$newSheet.Range("D1").Value = "Invalid E-Mail Address Text"
$newSheet.Range("D2").Value = "Please enter a valid e-mail address."

$newSheet.Range("E2").Value = "Please check your E-mail inbox for instructions to access your account."
$newSheet.Range("E1").Value = "Password Sent Header Text"

# Column widths to match the authored sheet.
$newSheet.Columns.Item(1).ColumnWidth = 26
$newSheet.Columns.Item(2).ColumnWidth = 25.7109375
$newSheet.Columns.Item(3).ColumnWidth = 26.140625
$newSheet.Columns.Item(4).ColumnWidth = 33
$newSheet.Columns.Item(5).ColumnWidth = 64.42578125

# Make the newly added sheet the active / selected tab.
$newSheet.Activate()
